$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Priming writes: establish the exact shared-string insertion order
$ws.Range("A5").Value = 'used'
$ws.Range("D4").Value = '1.3 VELOZ M/T'
$ws.Range("G4").Value = 'BINTARO'
$ws.Range("H4").Value = 'Pribadi'
$ws.Range("G5").Value = 'PLUIT'
$ws.Range("J3").Value = 'Personal'
$ws.Range("B5").Value = 'DAIHATSU'
$ws.Range("C5").Value = 'AYLA'
$ws.Range("G7").Value = 'HARAPAN INDAH'
$ws.Range("D5").Value = '1.0 D+M/T NEW'
$ws.Range("I3").Value = 'Finance'

# Remaining cells (reuse existing shared strings / literal numbers)
# Row 3
$ws.Range("A3").Value = 'new'
$ws.Range("B3").Value = 'TOYOTA'
$ws.Range("C3").Value = 'ALL NEW AVANZA'
$ws.Range("D3").Value = '1.3 E A/T'
$ws.Range("E3").Value = 2019
$ws.Range("F3").Value = 'DKI JAKARTA, BANTEN, JAWA BARAT'
$ws.Range("G3").Value = 'KARAWACI'
$ws.Range("H3").Value = 'Usaha'
$ws.Range("M3").Value = 'Yes'
$ws.Range("N3").Value = 'Yes'
$ws.Range("P3").Value = 36
$ws.Range("Q3").Value = 'passed'
# Row 4
$ws.Range("A4").Value = 'new'
$ws.Range("B4").Value = 'TOYOTA'
$ws.Range("C4").Value = 'ALL NEW AVANZA'
$ws.Range("E4").Value = 2019
$ws.Range("F4").Value = 'DKI JAKARTA, BANTEN, JAWA BARAT'
$ws.Range("J4").Value = 'Personal'
$ws.Range("L4").Value = 'Yes'
$ws.Range("M4").Value = 'Yes'
$ws.Range("N4").Value = 'Yes'
$ws.Range("P4").Value = 36
$ws.Range("Q4").Value = 'passed'
# Row 5
$ws.Range("E5").Value = 2018
$ws.Range("F5").Value = 'DKI JAKARTA, BANTEN, JAWA BARAT'
$ws.Range("H5").Value = 'Usaha'
$ws.Range("I5").Value = 'Pertanian'
$ws.Range("J5").Value = 'Company'
$ws.Range("O5").Value = 'Yes'
$ws.Range("P5").Value = 24
$ws.Range("Q5").Value = 'passed'
# Row 6
$ws.Range("A6").Value = 'used'
$ws.Range("B6").Value = 'DAIHATSU'
$ws.Range("C6").Value = 'AYLA'
$ws.Range("D6").Value = '1.0 D+M/T NEW'
$ws.Range("E6").Value = 2018
$ws.Range("F6").Value = 'DKI JAKARTA, BANTEN, JAWA BARAT'
$ws.Range("G6").Value = 'PLUIT'
$ws.Range("H6").Value = 'Usaha'
$ws.Range("I6").Value = 'Finance'
$ws.Range("J6").Value = 'Personal'
$ws.Range("O6").Value = 'Yes'
$ws.Range("P6").Value = 24
$ws.Range("Q6").Value = 'passed'
# Row 7
$ws.Range("A7").Value = 'used'
$ws.Range("B7").Value = 'DAIHATSU'
$ws.Range("C7").Value = 'AYLA'
$ws.Range("D7").Value = '1.0 D+M/T NEW'
$ws.Range("E7").Value = 2018
$ws.Range("F7").Value = 'DKI JAKARTA, BANTEN, JAWA BARAT'
$ws.Range("H7").Value = 'Pribadi'
$ws.Range("J7").Value = 'Personal'
$ws.Range("L7").Value = 'Yes'
$ws.Range("M7").Value = 'Yes'
$ws.Range("P7").Value = 60
$ws.Range("Q7").Value = 'passed'

# Refresh column widths for the new data range (approximate best-fit; exact
# sub-unit fractions aren't reproducible in this runtime's width model)
$ws.Columns.Item(1).ColumnWidth = 10.592447916666666
$ws.Columns.Item(2).ColumnWidth = 9.166666666666666
$ws.Columns.Item(3).ColumnWidth = 16.022135416666668
$ws.Columns.Item(4).ColumnWidth = 15.307291666666666
$ws.Columns.Item(5).ColumnWidth = 5.307291666666667
$ws.Columns.Item(6).ColumnWidth = 33.022135416666664
$ws.Columns.Item(7).ColumnWidth = 15.451822916666666
$ws.Columns.Item(8).ColumnWidth = 6.451822916666667
$ws.Columns.Item(9).ColumnWidth = 11.451822916666666
$ws.Columns.Item(10).ColumnWidth = 12.877604166666666
$ws.Columns.Item(11).ColumnWidth = 10.877604166666666
$ws.Columns.Item(12).ColumnWidth = 10.307291666666666
$ws.Columns.Item(13).ColumnWidth = 10.592447916666666
$ws.Columns.Item(14).ColumnWidth = 11.307291666666666
$ws.Columns.Item(15).ColumnWidth = 11.592447916666666
$ws.Columns.Item(16).ColumnWidth = 13.877604166666666
$ws.Columns.Item(17).ColumnWidth = 8.592447916666666

# Update selection to D4 (matches target sheetView/selection)
$ws.Range("D4").Select()
